# Convert from arrow batch to serialisation
# ---------------------------------------------------------------------------
# This script reproduces, on the "Consolidated" worksheet, the addition of a
# new asset row (Lorry/Tractor, inserted at row 19) plus nine brand-new asset
# rows appended at the bottom of the table (rows 39-47), and then promotes
# the A1:E47 range to a real Excel Table ("Table1") with an AutoFilter.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consolidated")

# --- 1. Insert a new row at 19 (shifts old rows 19-37 down to 20-38) -------
$ws.Rows("19:19").Insert()

$ws.Range("A19").Value = "Lorry/Tractor"
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = "Lorry/Tractor"
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = "Modern"

# --- 2. Append nine new asset rows (39-47) ----------------------------------
$ws.Range("A39").Value = "Swimming pool"
$ws.Range("B39").Value = 109
$ws.Range("C39").Value = "Swimming pool"
$ws.Range("D39").Value = 109
$ws.Range("E39").Value = 0

$ws.Range("A40").Value = "Air conditioner"
$ws.Range("B40").Value = 110
$ws.Range("C40").Value = "Appliance"
$ws.Range("D40").Value = 12
$ws.Range("E40").Formula = "=VLOOKUP(D40,AssetIdx,3,FALSE)"

$ws.Range("A41").Value = "Computer/Desktop/Laptop"
$ws.Range("B41").Value = 111
$ws.Range("C41").Value = "Appliance"
$ws.Range("D41").Value = 12
$ws.Range("E41").Formula = "=VLOOKUP(D41,AssetIdx,3,FALSE)"

$ws.Range("A42").Value = "Vacuum cleaner/Floor polisher"
$ws.Range("B42").Value = 112
$ws.Range("C42").Value = "Appliance"
$ws.Range("D42").Value = 12
$ws.Range("E42").Formula = "=VLOOKUP(D42,AssetIdx,3,FALSE)"

$ws.Range("A43").Value = "Dish washing machine"
$ws.Range("B43").Value = 113
$ws.Range("C43").Value = "Appliance"
$ws.Range("D43").Value = 12
$ws.Range("E43").Formula = "=VLOOKUP(D43,AssetIdx,3,FALSE)"

$ws.Range("A44").Value = "Tumble dryer"
$ws.Range("B44").Value = 114
$ws.Range("C44").Value = "Appliance"
$ws.Range("D44").Value = 12
$ws.Range("E44").Formula = "=VLOOKUP(D44,AssetIdx,3,FALSE)"

$ws.Range("A45").Value = "Home security service"
$ws.Range("B45").Value = 115
$ws.Range("C45").Value = "Home security service"
$ws.Range("D45").Value = 115
$ws.Range("E45").Value = 0

$ws.Range("A46").Value = "Geyser (providing hot water)"
$ws.Range("B46").Value = 116
$ws.Range("C46").Value = "Appliance"
$ws.Range("D46").Value = 12
$ws.Range("E46").Formula = "=VLOOKUP(D46,AssetIdx,3,FALSE)"

$ws.Range("A47").Value = "Pigs"
$ws.Range("B47").Value = 120
$ws.Range("C47").Value = "OtherLivestock"
$ws.Range("D47").Value = 28
$ws.Range("E47").Formula = "=VLOOKUP(D47,AssetIdx,3,FALSE)"

# --- 3. Re-apply the (unchanged) two-key sort over the old data block ------
#     so the worksheet's recorded sortState range grows from E37 to E38.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B38"))
$ws.Sort.SortFields.Add($ws.Range("D2:D38"))
$ws.Sort.SetRange($ws.Range("A2:E38"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- 4. Promote A1:E47 to an Excel Table (ListObject) -----------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:E47"), 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"

# --- 5. Column widths (best effort; engine quantizes to 1/6 char units) ----
$ws.Columns.Item(1).ColumnWidth = 26.264322916666668
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(5).ColumnWidth = 9.565104166666666
$ws.Columns.Item(7).ColumnWidth = 31.565104166666668

# --- 6. Selection matches the post-edit UI state ----------------------------
$ws.Range("A39:E47").Select()
